$wb = $excel.ActiveWorkbook

# Sheet 1: ALC
$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(17, 8).Value = 181246.33  # H17: 186275.94 -> 181246.33
$ws.Cells.Item(17, 10).Value = 186253.17  # J17: 191569.55 -> 186253.17
$ws.Cells.Item(17, 12).Value = 558759.51  # L17: 574708.6499999999 -> 558759.51
$ws.Cells.Item(17, 14).Value = -559095.51  # N17: -575044.6499999999 -> -559095.51
$ws.Cells.Item(28, 8).Value = 523.7895  # H28: 568 -> 523.7895
$ws.Cells.Item(28, 9).Value = 539.8125  # I28: 575.13336 -> 539.8125
$ws.Cells.Item(28, 10).Value = 438.33334  # J28: 514.5 -> 438.33334
$ws.Cells.Item(28, 11).Value = 539.8125  # K28: 575.13336 -> 539.8125
$ws.Cells.Item(28, 12).Value = 438.33334  # L28: 514.5 -> 438.33334
$ws.Cells.Item(28, 13).Value = -54.8125  # M28: -90.13336000000004 -> -54.8125
$ws.Cells.Item(28, 14).Value = -1408.33334  # N28: -1484.5 -> -1408.33334
$ws.Cells.Item(113, 8).Value = 4657.2  # H113: 4833.3335 -> 4657.2
$ws.Cells.Item(113, 10).Value = 4595.3335  # J113: 5000 -> 4595.3335
$ws.Cells.Item(113, 12).Value = 4595.3335  # L113: 5000 -> 4595.3335
$ws.Cells.Item(113, 14).Value = -11103.3335  # N113: -11508 -> -11103.3335
$ws.Cells.Item(132, 8).Value = 1537.0476  # H132: 1685.2222 -> 1537.0476
$ws.Cells.Item(132, 9).Value = 1434  # I132: 1575.6 -> 1434
$ws.Cells.Item(132, 10).Value = 1975  # J132: 2233.3333 -> 1975
$ws.Cells.Item(132, 11).Value = 4302  # K132: 4726.799999999999 -> 4302
$ws.Cells.Item(132, 12).Value = 5925  # L132: 6699.999899999999 -> 5925
$ws.Cells.Item(132, 13).Value = -1772  # M132: -2196.799999999999 -> -1772
$ws.Cells.Item(132, 14).Value = -10985  # N132: -11759.9999 -> -10985

# Sheet 2: ARM
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(12, 8).Value = 2000  # H12: 0 -> 2000
$ws.Cells.Item(12, 9).Value = 2000  # I12: 0 -> 2000
$ws.Cells.Item(12, 11).Value = 2000  # K12: 0 -> 2000
$ws.Cells.Item(12, 13).Value = -1827  # M12: None -> -1827
$ws.Cells.Item(16, 8).Value = 4998  # H16: 0 -> 4998
$ws.Cells.Item(16, 9).Value = 4998  # I16: 0 -> 4998
$ws.Cells.Item(16, 11).Value = 4998  # K16: 0 -> 4998
$ws.Cells.Item(16, 13).Value = -4711  # M16: None -> -4711
$ws.Cells.Item(36, 8).Value = 7598.6  # H36: 6815.5 -> 7598.6
$ws.Cells.Item(36, 9).Value = 7598.6  # I36: 6815.5 -> 7598.6
$ws.Cells.Item(36, 11).Value = 7598.6  # K36: 6815.5 -> 7598.6
$ws.Cells.Item(36, 13).Value = -7252.6  # M36: -6469.5 -> -7252.6
$ws.Cells.Item(45, 8).Value = 92918.91  # H45: 127013.625 -> 92918.91
$ws.Cells.Item(45, 9).Value = 126875  # I45: 201800.2 -> 126875
$ws.Cells.Item(45, 11).Value = 126875  # K45: 201800.2 -> 126875
$ws.Cells.Item(45, 13).Value = -126498  # M45: -201423.2 -> -126498
$ws.Cells.Item(61, 8).Value = 27882.512  # H61: 26591.104 -> 27882.512
$ws.Cells.Item(61, 9).Value = 39689.55  # I61: 36645.53 -> 39689.55
$ws.Cells.Item(61, 11).Value = 39689.55  # K61: 36645.53 -> 39689.55
$ws.Cells.Item(61, 13).Value = -39477.55  # M61: -36433.53 -> -39477.55
$ws.Cells.Item(74, 8).Value = 22933.688  # H74: 23898.195 -> 22933.688
$ws.Cells.Item(74, 9).Value = 2071.3242  # I74: 2101.0833 -> 2071.3242
$ws.Cells.Item(74, 10).Value = 93107.09  # J74: 102367.8 -> 93107.09
$ws.Cells.Item(74, 11).Value = 2071.3242  # K74: 2101.0833 -> 2071.3242
$ws.Cells.Item(74, 12).Value = 93107.09  # L74: 102367.8 -> 93107.09
$ws.Cells.Item(74, 13).Value = -1197.3242  # M74: -1227.0833 -> -1197.3242
$ws.Cells.Item(74, 14).Value = -94855.09  # N74: -104115.8 -> -94855.09
$ws.Cells.Item(77, 8).Value = 22933.688  # H77: 23898.195 -> 22933.688
$ws.Cells.Item(77, 9).Value = 2071.3242  # I77: 2101.0833 -> 2071.3242
$ws.Cells.Item(77, 10).Value = 93107.09  # J77: 102367.8 -> 93107.09
$ws.Cells.Item(77, 11).Value = 10356.621  # K77: 10505.4165 -> 10356.621
$ws.Cells.Item(77, 12).Value = 465535.45  # L77: 511839 -> 465535.45
$ws.Cells.Item(77, 13).Value = -5988.620999999999  # M77: -6137.416499999999 -> -5988.620999999999
$ws.Cells.Item(77, 14).Value = -474271.45  # N77: -520575 -> -474271.45
$ws.Cells.Item(122, 8).Value = 1632.1765  # H122: 1553.4286 -> 1632.1765
$ws.Cells.Item(122, 9).Value = 803.3570999999999  # I122: 854 -> 803.3570999999999
$ws.Cells.Item(122, 10).Value = 5500  # J122: 5750 -> 5500
$ws.Cells.Item(122, 11).Value = 2410.0713  # K122: 2562 -> 2410.0713
$ws.Cells.Item(122, 12).Value = 16500  # L122: 17250 -> 16500
$ws.Cells.Item(122, 13).Value = 39.92870000000039  # M122: -112 -> 39.92870000000039
$ws.Cells.Item(122, 14).Value = -21400  # N122: -22150 -> -21400
$ws.Cells.Item(136, 8).Value = 27882.512  # H136: 26591.104 -> 27882.512
$ws.Cells.Item(136, 9).Value = 39689.55  # I136: 36645.53 -> 39689.55
$ws.Cells.Item(136, 11).Value = 119068.65  # K136: 109936.59 -> 119068.65
$ws.Cells.Item(136, 13).Value = -116518.65  # M136: -107386.59 -> -116518.65

# Sheet 3: BSM
$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(19, 8).Value = 38000  # H19: 31100 -> 38000
$ws.Cells.Item(19, 9).Value = 38000  # I19: 31500 -> 38000
$ws.Cells.Item(19, 10).Value = 0  # J19: 29500 -> 0
$ws.Cells.Item(19, 11).Value = 38000  # K19: 31500 -> 38000
$ws.Cells.Item(19, 12).Value = 0  # L19: 29500 -> 0
$ws.Cells.Item(19, 13).Value = -37827  # M19: -31327 -> -37827
$ws.Cells.Item(19, 14).ClearContents()  # was N19=-29846
$ws.Cells.Item(134, 8).Value = 1500.7551  # H134: 1498.3673 -> 1500.7551
$ws.Cells.Item(134, 9).Value = 1339.4857  # I134: 1336.1428 -> 1339.4857
$ws.Cells.Item(134, 11).Value = 4018.4571  # K134: 4008.4284 -> 4018.4571
$ws.Cells.Item(134, 13).Value = -1483.4571  # M134: -1473.4284 -> -1483.4571

# Sheet 4: CRP
$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(31, 8).Value = 3618.75  # H31: 3681.25 -> 3618.75
$ws.Cells.Item(31, 9).Value = 2605  # I31: 2709.1667 -> 2605
$ws.Cells.Item(31, 11).Value = 2605  # K31: 2709.1667 -> 2605
$ws.Cells.Item(31, 13).Value = -2310  # M31: -2414.1667 -> -2310
$ws.Cells.Item(34, 8).Value = 3618.75  # H34: 3681.25 -> 3618.75
$ws.Cells.Item(34, 9).Value = 2605  # I34: 2709.1667 -> 2605
$ws.Cells.Item(34, 11).Value = 2605  # K34: 2709.1667 -> 2605
$ws.Cells.Item(34, 13).Value = -2403  # M34: -2507.1667 -> -2403
$ws.Cells.Item(50, 8).Value = 0  # H50: 49999 -> 0
$ws.Cells.Item(50, 10).Value = 0  # J50: 49999 -> 0
$ws.Cells.Item(50, 12).Value = 0  # L50: 49999 -> 0
$ws.Cells.Item(50, 14).ClearContents()  # was N50=-51249
$ws.Cells.Item(99, 8).Value = 4001082.2  # H99: 10001150 -> 4001082.2
$ws.Cells.Item(99, 9).Value = 5001103  # I99: 10001150 -> 5001103
$ws.Cells.Item(99, 10).Value = 1000  # J99: 0 -> 1000
$ws.Cells.Item(99, 11).Value = 5001103  # K99: 10001150 -> 5001103
$ws.Cells.Item(99, 12).Value = 1000  # L99: 0 -> 1000
$ws.Cells.Item(99, 13).Value = -4999605  # M99: -9999652 -> -4999605
$ws.Cells.Item(99, 14).Value = -3996  # N99: None -> -3996
$ws.Cells.Item(122, 8).Value = 0  # H122: 1000 -> 0
$ws.Cells.Item(122, 9).Value = 0  # I122: 1000 -> 0
$ws.Cells.Item(122, 11).Value = 0  # K122: 3000 -> 0
$ws.Cells.Item(122, 13).ClearContents()  # was M122=-550
$ws.Cells.Item(126, 8).Value = 4001082.2  # H126: 10001150 -> 4001082.2
$ws.Cells.Item(126, 9).Value = 5001103  # I126: 10001150 -> 5001103
$ws.Cells.Item(126, 10).Value = 1000  # J126: 0 -> 1000
$ws.Cells.Item(126, 11).Value = 15003309  # K126: 30003450 -> 15003309
$ws.Cells.Item(126, 12).Value = 3000  # L126: 0 -> 3000
$ws.Cells.Item(126, 13).Value = -15000839  # M126: -30000980 -> -15000839
$ws.Cells.Item(126, 14).Value = -7940  # N126: None -> -7940
$ws.Cells.Item(141, 8).Value = 273241  # H141: 273241.84 -> 273241
$ws.Cells.Item(141, 9).Value = 77121  # I141: 0 -> 77121
$ws.Cells.Item(141, 10).Value = 312465  # J141: 273241.84 -> 312465
$ws.Cells.Item(141, 11).Value = 77121  # K141: 0 -> 77121
$ws.Cells.Item(141, 12).Value = 312465  # L141: 273241.84 -> 312465
$ws.Cells.Item(141, 14).Value = -322825  # N141: -283601.84 -> -322825
$ws.Cells.Item(141, 13).Value = -71941  # M141: None -> -71941

# Sheet 5: CUL
$ws = $wb.Worksheets.Item(5)
$ws.Cells.Item(33, 8).Value = 20100246  # H33: 18273772 -> 20100246
$ws.Cells.Item(33, 10).Value = 22333570  # J33: 20101116 -> 22333570
$ws.Cells.Item(33, 12).Value = 134001420  # L33: 120606696 -> 134001420
$ws.Cells.Item(33, 14).Value = -134001986  # N33: -120607262 -> -134001986
$ws.Cells.Item(121, 8).Value = 10269.85  # H121: 9818.904 -> 10269.85
$ws.Cells.Item(121, 10).Value = 12392.1875  # J121: 11710.294 -> 12392.1875
$ws.Cells.Item(121, 12).Value = 37176.5625  # L121: 35130.882 -> 37176.5625
$ws.Cells.Item(121, 14).Value = -39796.5625  # N121: -37750.882 -> -39796.5625

# Sheet 6: GSM
$ws = $wb.Worksheets.Item(6)
$ws.Cells.Item(36, 8).Value = 2800  # H36: 1860 -> 2800
$ws.Cells.Item(36, 10).Value = 2933.3333  # J36: 650 -> 2933.3333
$ws.Cells.Item(36, 12).Value = 2933.3333  # L36: 650 -> 2933.3333
$ws.Cells.Item(36, 14).Value = -3903.3333  # N36: -1620 -> -3903.3333
$ws.Cells.Item(80, 8).Value = 20021108  # H80: 22751096 -> 20021108
$ws.Cells.Item(80, 9).Value = 21513.867  # I80: 22965.357 -> 21513.867
$ws.Cells.Item(80, 10).Value = 50020500  # J80: 62525324 -> 50020500
$ws.Cells.Item(80, 11).Value = 21513.867  # K80: 22965.357 -> 21513.867
$ws.Cells.Item(80, 12).Value = 50020500  # L80: 62525324 -> 50020500
$ws.Cells.Item(80, 13).Value = -20515.867  # M80: -21967.357 -> -20515.867
$ws.Cells.Item(80, 14).Value = -50022496  # N80: -62527320 -> -50022496
$ws.Cells.Item(83, 8).Value = 20021108  # H83: 22751096 -> 20021108
$ws.Cells.Item(83, 9).Value = 21513.867  # I83: 22965.357 -> 21513.867
$ws.Cells.Item(83, 10).Value = 50020500  # J83: 62525324 -> 50020500
$ws.Cells.Item(83, 11).Value = 107569.335  # K83: 114826.785 -> 107569.335
$ws.Cells.Item(83, 12).Value = 250102500  # L83: 312626620 -> 250102500
$ws.Cells.Item(83, 13).Value = -102577.335  # M83: -109834.785 -> -102577.335
$ws.Cells.Item(83, 14).Value = -250112484  # N83: -312636604 -> -250112484
$ws.Cells.Item(102, 8).Value = 16668102  # H102: 17858608 -> 16668102
$ws.Cells.Item(102, 9).Value = 20834664  # I102: 21740494 -> 20834664
$ws.Cells.Item(102, 10).Value = 1857.6666  # J102: 1929.2 -> 1857.6666
$ws.Cells.Item(102, 11).Value = 20834664  # K102: 21740494 -> 20834664
$ws.Cells.Item(102, 12).Value = 1857.6666  # L102: 1929.2 -> 1857.6666
$ws.Cells.Item(102, 13).Value = -20833042  # M102: -21738872 -> -20833042
$ws.Cells.Item(102, 14).Value = -5101.6666  # N102: -5173.2 -> -5101.6666
$ws.Cells.Item(107, 8).Value = 143949.42  # H107: 126005.75 -> 143949.42
$ws.Cells.Item(107, 10).Value = 1594.5  # J107: 1355.6 -> 1594.5
$ws.Cells.Item(107, 12).Value = 1594.5  # L107: 1355.6 -> 1594.5
$ws.Cells.Item(107, 14).Value = -5434.5  # N107: -5195.6 -> -5434.5
$ws.Cells.Item(122, 8).Value = 2379  # H122: 2049.75 -> 2379
$ws.Cells.Item(122, 9).Value = 2230.8572  # I122: 1902.1666 -> 2230.8572
$ws.Cells.Item(122, 10).Value = 2793.8  # J122: 2492.5 -> 2793.8
$ws.Cells.Item(122, 11).Value = 6692.571599999999  # K122: 5706.4998 -> 6692.571599999999
$ws.Cells.Item(122, 12).Value = 8381.400000000001  # L122: 7477.5 -> 8381.400000000001
$ws.Cells.Item(122, 13).Value = -4242.571599999999  # M122: -3256.4998 -> -4242.571599999999
$ws.Cells.Item(122, 14).Value = -13281.4  # N122: -12377.5 -> -13281.4

# Sheet 7: LTW
$ws = $wb.Worksheets.Item(7)
$ws.Cells.Item(12, 8).Value = 2000  # H12: 5000 -> 2000
$ws.Cells.Item(12, 9).Value = 0  # I12: 5000 -> 0
$ws.Cells.Item(12, 10).Value = 2000  # J12: 0 -> 2000
$ws.Cells.Item(12, 11).Value = 0  # K12: 5000 -> 0
$ws.Cells.Item(12, 12).Value = 2000  # L12: 0 -> 2000
$ws.Cells.Item(12, 13).ClearContents()  # was M12=-4830
$ws.Cells.Item(12, 14).Value = -2340  # N12: None -> -2340
$ws.Cells.Item(46, 8).Value = 7471.28  # H46: 6754.517 -> 7471.28
$ws.Cells.Item(46, 9).Value = 22899.6  # I46: 15062.125 -> 22899.6
$ws.Cells.Item(46, 10).Value = 3614.2  # J46: 3589.7144 -> 3614.2
$ws.Cells.Item(46, 11).Value = 22899.6  # K46: 15062.125 -> 22899.6
$ws.Cells.Item(46, 12).Value = 3614.2  # L46: 3589.7144 -> 3614.2
$ws.Cells.Item(46, 13).Value = -22711.6  # M46: -14874.125 -> -22711.6
$ws.Cells.Item(46, 14).Value = -3990.2  # N46: -3965.7144 -> -3990.2
$ws.Cells.Item(61, 8).Value = 18681.666  # H61: 17511 -> 18681.666
$ws.Cells.Item(61, 9).Value = 18022.5  # I61: 16681.334 -> 18022.5
$ws.Cells.Item(61, 11).Value = 18022.5  # K61: 16681.334 -> 18022.5
$ws.Cells.Item(61, 13).Value = -17820.5  # M61: -16479.334 -> -17820.5
$ws.Cells.Item(113, 8).Value = 18681.666  # H113: 17511 -> 18681.666
$ws.Cells.Item(113, 9).Value = 18022.5  # I113: 16681.334 -> 18022.5
$ws.Cells.Item(113, 11).Value = 18022.5  # K113: 16681.334 -> 18022.5
$ws.Cells.Item(113, 13).Value = -15852.5  # M113: -14511.334 -> -15852.5
$ws.Cells.Item(122, 8).Value = 2787.0908  # H122: 2830.4375 -> 2787.0908
$ws.Cells.Item(122, 9).Value = 2573.4  # I122: 2622.2917 -> 2573.4
$ws.Cells.Item(122, 11).Value = 7720.200000000001  # K122: 7866.875100000001 -> 7720.200000000001
$ws.Cells.Item(122, 13).Value = -5270.200000000001  # M122: -5416.875100000001 -> -5270.200000000001
$ws.Cells.Item(132, 8).Value = 3550.389  # H132: 3747.1765 -> 3550.389
$ws.Cells.Item(132, 9).Value = 2672.923  # I132: 2749.3845 -> 2672.923
$ws.Cells.Item(132, 10).Value = 5831.8  # J132: 6990 -> 5831.8
$ws.Cells.Item(132, 11).Value = 8018.768999999999  # K132: 8248.1535 -> 8018.768999999999
$ws.Cells.Item(132, 12).Value = 17495.4  # L132: 20970 -> 17495.4
$ws.Cells.Item(132, 13).Value = -5488.768999999999  # M132: -5718.1535 -> -5488.768999999999
$ws.Cells.Item(132, 14).Value = -22555.4  # N132: -26030 -> -22555.4

# Sheet 8: WVR
$ws = $wb.Worksheets.Item(8)
$ws.Cells.Item(12, 8).Value = 5000  # H12: 4937.3335 -> 5000
$ws.Cells.Item(12, 9).Value = 5000  # I12: 4937.3335 -> 5000
$ws.Cells.Item(12, 11).Value = 5000  # K12: 4937.3335 -> 5000
$ws.Cells.Item(12, 13).Value = -4858  # M12: -4795.3335 -> -4858
$ws.Cells.Item(107, 8).Value = 1431097.6  # H107: 1506303.4 -> 1431097.6
$ws.Cells.Item(107, 9).Value = 2346  # I107: 2372.1667 -> 2346
$ws.Cells.Item(107, 11).Value = 7038  # K107: 7116.500100000001 -> 7038
$ws.Cells.Item(107, 13).Value = -5118  # M107: -5196.500100000001 -> -5118
$ws.Cells.Item(122, 8).Value = 5834.375  # H122: 6132.4 -> 5834.375
$ws.Cells.Item(122, 9).Value = 6133.2  # I122: 6165.5 -> 6133.2
$ws.Cells.Item(122, 10).Value = 5336.3335  # J122: 6000 -> 5336.3335
$ws.Cells.Item(122, 11).Value = 18399.6  # K122: 18496.5 -> 18399.6
$ws.Cells.Item(122, 12).Value = 16009.0005  # L122: 18000 -> 16009.0005
$ws.Cells.Item(122, 13).Value = -15949.6  # M122: -16046.5 -> -15949.6
$ws.Cells.Item(122, 14).Value = -20909.0005  # N122: -22900 -> -20909.0005
$ws.Cells.Item(136, 8).Value = 1742  # H136: 1834.3334 -> 1742
$ws.Cells.Item(136, 9).Value = 1090.4  # I136: 1201.2 -> 1090.4
$ws.Cells.Item(136, 11).Value = 3271.2  # K136: 3603.6 -> 3271.2
$ws.Cells.Item(136, 13).Value = -721.2000000000003  # M136: -1053.6 -> -721.2000000000003
